# Update gh-pages output values for the "展览" and "全部类型" sheets.
# Both sheets carry the same underlying event list, so the same edits are
# applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column F ("想去人数") bumps for several rows.
    $ws.Range("F2").Value  = 1057
    $ws.Range("F3").Value  = 748
    $ws.Range("F8").Value  = 1684
    $ws.Range("F9").Value  = 6381
    $ws.Range("F11").Value = 365
    $ws.Range("F12").Value = 299
    $ws.Range("F13").Value = 96
    $ws.Range("F14").Value = 375
    $ws.Range("F16").Value = 6350
    $ws.Range("F17").Value = 272
    $ws.Range("F18").Value = 1281
    $ws.Range("F22").Value = 105
    $ws.Range("F23").Value = 272
    $ws.Range("F33").Value = 46

    # Row 20 event got cancelled: name suffix, lower F count, and the
    # lowest-ticket-price column becomes a text marker instead of a number.
    $ws.Range("C20").Value = "江西·2024南昌玛雅《次元之芯》主题动漫嘉年华（取消）"
    $ws.Range("F20").Value = 116
    $ws.Range("G20").Value = "不可售"
}
